function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws.Range('D2') '29.912.04'
$ws.Range('E2').Value = '  -0.96%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.898.32'
$ws.Range('E3').Value = '  -0.60%  '

# Row 4
$ws.Range('E4').Value = '  -0.38%  '

# Row 5
Set-TextValue $ws.Range('D5') '0.7577'
$ws.Range('E5').Value = '  +3.62%  '

# Row 6
Set-TextValue $ws.Range('D6') '239.83'
$ws.Range('E6').Value = '  -1.54%  '

# Row 7
Set-TextValue $ws.Range('D7') '1.000'
$ws.Range('E7').Value = '  -0.30%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.3048'
$ws.Range('E8').Value = '  -2.50%  '

# Row 9
Set-TextValue $ws.Range('D9') '25.43'
$ws.Range('E9').Value = '  -5.07%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.06844'
$ws.Range('E10').Value = '  -1.09%  '

# Row 11
$ws.Range('E11').Value = '  -0.10%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.7444'
$ws.Range('E12').Value = '  -4.25%  '

# Row 13
Set-TextValue $ws.Range('D13') '1.898.87'
$ws.Range('E13').Value = '  -0.33%  '

# Row 14
Set-TextValue $ws.Range('D14') '5.179'
$ws.Range('E14').Value = '  -1.49%  '

# Row 15
Set-TextValue $ws.Range('D15') '91.17'
$ws.Range('E15').Value = '  -0.23%  '

# Row 16
Set-TextValue $ws.Range('D16') '29.920.43'
$ws.Range('E16').Value = '  -0.68%  '

# Row 17
Set-TextValue $ws.Range('D17') '13.94'
$ws.Range('E17').Value = '  -1.92%  '

# Row 18
Set-TextValue $ws.Range('D18') '5.958'
$ws.Range('E18').Value = '  +1.93%  '

# Row 19
Set-TextValue $ws.Range('D19') '243.00'
$ws.Range('E19').Value = '  +0.90%  '

# Row 20
Set-TextValue $ws.Range('D20') '0.000007712'
$ws.Range('E20').Value = '  -1.03%  '

# Row 21
$ws.Range('E21').Value = '  -0.22%  '

# Row 22
$ws.Range('E22').Value = '  -0.35%  '

# Row 23
Set-TextValue $ws.Range('D23') '6.922'
$ws.Range('E23').Value = '  +2.86%  '

# Row 24
Set-TextValue $ws.Range('D24') '9.244'
$ws.Range('E24').Value = '  -1.73%  '

# Row 25
Set-TextValue $ws.Range('D25') '165.43'
$ws.Range('E25').Value = '  -0.17%  '

# Row 26
Set-TextValue $ws.Range('D26') '18.71'
$ws.Range('E26').Value = '  -1.61%  '

# Row 27
Set-TextValue $ws.Range('D27') '0.1271'
$ws.Range('E27').Value = '  +0.22%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.018'
$ws.Range('E28').Value = '  -3.39%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.385'
$ws.Range('E29').Value = '  +2.13%  '

# Row 30
Set-TextValue $ws.Range('D30') '1.515'
$ws.Range('E30').Value = '  -2.06%  '

# Row 31
Set-TextValue $ws.Range('D31') '4.273'
$ws.Range('E31').Value = '  -0.73%  '

# Row 32
Set-TextValue $ws.Range('D32') '4.019'
$ws.Range('E32').Value = '  -1.59%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.05301'
$ws.Range('E33').Value = '  +2.70%  '

# Row 34
Set-TextValue $ws.Range('D34') '1.245'
$ws.Range('E34').Value = '  -2.87%  '

# Row 35
Set-TextValue $ws.Range('D35') '0.7230'
$ws.Range('E35').Value = '  -2.85%  '

# Row 36
Set-TextValue $ws.Range('D36') '2.715'
$ws.Range('E36').Value = '  -1.57%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.01910'
$ws.Range('E37').Value = '  -1.30%  '

# Row 38
$ws.Range('E38').Value = '  -0.30%  '

# Row 39
Set-TextValue $ws.Range('D39') '6.151'
$ws.Range('E39').Value = '  -3.33%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.4398'
$ws.Range('E40').Value = '  -1.08%  '

# Row 41
Set-TextValue $ws.Range('D41') '71.98'
$ws.Range('E41').Value = '  -3.60%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.9999'
$ws.Range('E42').Value = '  -0.20%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.879'
$ws.Range('E43').Value = '  -2.63%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.8252'
$ws.Range('E44').Value = '  -1.17%  '

# Row 45
Set-TextValue $ws.Range('D45') '100.70'
$ws.Range('E45').Value = '  -0.21%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D46') '9.824'
$ws.Range('E46').Value = '  +0.84%  '

# Row 47
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D47') '7.515'
$ws.Range('E47').Value = '  -1.07%  '

# Row 48
Set-TextValue $ws.Range('D48') '2.049.78'
$ws.Range('E48').Value = '  -0.21%  '

# Row 49
Set-TextValue $ws.Range('D49') '36.25'
$ws.Range('E49').Value = '  -3.46%  '

# Row 50
Set-TextValue $ws.Range('D50') '0.05963'
$ws.Range('E50').Value = '  -0.38%  '

# Row 51
Set-TextValue $ws.Range('D51') '1.474'
$ws.Range('E51').Value = '  +1.17%  '
